# Auto-generated files on 2025-08-31
# Update the HotStock_Top20 ranking table (columns A, B, C for rows 2-21)
# to reflect the latest hot-stock rankings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("通富微电", "通富微电", "岩山科技"),
    @("寒武纪-U", "长电科技", "华胜天成"),
    @("数据港",   "数据港",   "通富微电"),
    @("先导智能", "岩山科技", "长城军工"),
    @("岩山科技", "寒武纪-U", "领益智造"),
    @("工业富联", "国轩高科", "工业富联"),
    @("全志科技", "全志科技", "吉视传媒"),
    @("长电科技", "先导智能", "国轩高科"),
    @("中芯国际", "东方财富", "三维通信"),
    @("国轩高科", "中国稀土", "寒武纪"),
    @("三维通信", "寒武纪-U", "长电科技"),
    @("吉视传媒", "三维通信", "北方稀土"),
    @("华胜天成", "吉视传媒", "数据港"),
    @("长城军工", "华工科技", "先导智能"),
    @("中国稀土", "旋极信息", "拓维信息"),
    @("领益智造", "长城军工", "中国稀土"),
    @("拓维信息", "华胜天成", "天融信"),
    @("北方稀土", "领益智造", "利欧股份"),
    @("东方财富", "中芯国际", "中芯国际"),
    @("旋极信息", "拓普集团", "歌尔股份")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
